# Week 16 log + new player (T.Lewis) roster addition, applied to both the
# "Rushing" and "Receiving" sheets.
#
# The new player "T.Lewis" is inserted as a new column right before the
# existing "M.Gesicki" column (i.e. immediately after "K.Merritt"), so the
# player order on each sheet stays: ... K.Merritt, T.Lewis, M.Gesicki,
# A.Shaheen, H.Long, D.Smythe. Row 2 keeps logging "n" for every player,
# including the newly added one.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Find the column whose header (row 1) is "M.Gesicki" - the new player
    # column gets inserted right before it.
    $headerRange = $ws.Range("A1:Z1")
    $found = $headerRange.Find("M.Gesicki")
    $col = $found.Column

    # Insert a new column at that position; existing M.Gesicki..D.Smythe
    # columns (and everything after) shift one column to the right,
    # carrying their formatting/values with them.
    $ws.Columns($col).Insert()

    $newCell = $ws.Cells.Item(1, $col)
    $newCell.Value = "T.Lewis"

    $dataCell = $ws.Cells.Item(2, $col)
    $dataCell.Value = "n"
}
